$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (sheet1.xml) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1196
$ws.Range("J17").Value = 1196
$ws.Range("L17").Value = 3588
$ws.Range("N17").Value = -3924
$ws.Range("H29").Value = 3321.8333
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H43").Value = 5166.3335
$ws.Range("J43").Value = 6749.5
$ws.Range("L43").Value = 6749.5
$ws.Range("N43").Value = -6887.5
$ws.Range("H92").Value = 110.625
$ws.Range("I92").Value = 123.333336
$ws.Range("J92").Value = 72.5
$ws.Range("K92").Value = 123.333336
$ws.Range("L92").Value = 72.5
$ws.Range("M92").Value = 1124.666664
$ws.Range("N92").Value = -2568.5
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800
$ws.Range("H131").Value = 2500
$ws.Range("I131").Value = 2500
$ws.Range("K131").Value = 7500
$ws.Range("M131").Value = -2460
$ws.Range("H135").Value = 1097
$ws.Range("I135").Value = 1163.1666
$ws.Range("J135").Value = 700
$ws.Range("K135").Value = 10468.4994
$ws.Range("L135").Value = 6300
$ws.Range("M135").Value = -7933.499400000001
$ws.Range("N135").Value = -11370

# ---- Sheet: ARM (sheet2.xml) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2775.2727
$ws.Range("J45").Value = 4433.6
$ws.Range("L45").Value = 4433.6
$ws.Range("N45").Value = -5187.6
$ws.Range("H58").Value = 41600
$ws.Range("I58").Value = 28000
$ws.Range("J58").Value = 62000
$ws.Range("K58").Value = 28000
$ws.Range("L58").Value = 62000
$ws.Range("M58").Value = -27570
$ws.Range("N58").Value = -62860
$ws.Range("H74").Value = 3034
$ws.Range("I74").Value = 3034
$ws.Range("K74").Value = 3034
$ws.Range("M74").Value = -2160
$ws.Range("H77").Value = 3034
$ws.Range("I77").Value = 3034
$ws.Range("K77").Value = 15170
$ws.Range("M77").Value = -10802
$ws.Range("H97").Value = 720
$ws.Range("I97").Value = 479.9
$ws.Range("J97").Value = 1320.25
$ws.Range("K97").Value = 479.9
$ws.Range("L97").Value = 1320.25
$ws.Range("M97").Value = 16.10000000000002
$ws.Range("N97").Value = -2312.25
$ws.Range("H102").Value = 5184.375
$ws.Range("I102").Value = 2458.6667
$ws.Range("K102").Value = 2458.6667
$ws.Range("M102").Value = -836.6667000000002
$ws.Range("H132").Value = 3728.724
$ws.Range("I132").Value = 3701.8572
$ws.Range("K132").Value = 11105.5716
$ws.Range("M132").Value = -8575.571599999999

# ---- Sheet: BSM (sheet3.xml) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2350
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2350
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").Value = 2350
$ws.Range("N80").Value = -4346
$ws.Range("H83").Value = 2350
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2350
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").Value = 11750
$ws.Range("N83").Value = -21734
$ws.Range("H86").Value = 7424.375
$ws.Range("I86").Value = 5795
$ws.Range("J86").Value = 7657.143
$ws.Range("K86").Value = 5795
$ws.Range("L86").Value = 7657.143
$ws.Range("M86").Value = -4672
$ws.Range("N86").Value = -9903.143
$ws.Range("H89").Value = 7424.375
$ws.Range("I89").Value = 5795
$ws.Range("J89").Value = 7657.143
$ws.Range("K89").Value = 28975
$ws.Range("L89").Value = 38285.715
$ws.Range("M89").Value = -23359
$ws.Range("N89").Value = -49517.715

# ---- Sheet: CRP (sheet4.xml) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 10000000
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H50").Value = 39770.75
$ws.Range("I50").Value = 20083
$ws.Range("J50").Value = 46333.332
$ws.Range("K50").Value = 20083
$ws.Range("L50").Value = 46333.332
$ws.Range("M50").Value = -19458
$ws.Range("N50").Value = -47583.332
$ws.Range("H51").Value = 98765
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H57").Value = 10000
$ws.Range("I57").Value = 10000
$ws.Range("K57").Value = 10000
$ws.Range("M57").Value = -9440
$ws.Range("H61").Value = 98765
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

# ---- Sheet: CUL (sheet5.xml) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 46588.766
$ws.Range("I4").Value = 66742.3
$ws.Range("J4").Value = 80.61539
$ws.Range("K4").Value = 200226.9
$ws.Range("L4").Value = 241.84617
$ws.Range("M4").Value = -200114.9
$ws.Range("N4").Value = -465.84617
$ws.Range("H60").Value = 1726.7858
$ws.Range("I60").Value = 123.333336
$ws.Range("J60").Value = 2164.0908
$ws.Range("K60").Value = 370.000008
$ws.Range("L60").Value = 6492.2724
$ws.Range("M60").Value = -119.000008
$ws.Range("N60").Value = -6994.2724
$ws.Range("H97").Value = 1363.75
$ws.Range("I97").Value = 450
$ws.Range("J97").Value = 1668.3334
$ws.Range("K97").Value = 1350
$ws.Range("L97").Value = 5005.0002
$ws.Range("M97").Value = -854
$ws.Range("N97").Value = -5997.0002
$ws.Range("H109").Value = 1953
$ws.Range("I109").Value = 1953
$ws.Range("K109").Value = 5859
$ws.Range("M109").Value = -4819

# ---- Sheet: GSM (sheet6.xml) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3816.2222
$ws.Range("J70").Value = 3886.5
$ws.Range("L70").Value = 3886.5
$ws.Range("N70").Value = -4426.5
$ws.Range("H73").Value = 3816.2222
$ws.Range("J73").Value = 3886.5
$ws.Range("L73").Value = 3886.5
$ws.Range("N73").Value = -5758.5
$ws.Range("H102").Value = 1986
$ws.Range("I102").Value = 1865.44
$ws.Range("K102").Value = 1865.44
$ws.Range("M102").Value = -243.4400000000001
$ws.Range("H113").Value = 3140.2
$ws.Range("I113").Value = 3140.2
$ws.Range("K113").Value = 3140.2
$ws.Range("M113").Value = -970.1999999999998
$ws.Range("H122").Value = 2582.4285
$ws.Range("I122").Value = 1860
$ws.Range("K122").Value = 5580
$ws.Range("M122").Value = -3130
$ws.Range("H132").Value = 3111.6155
$ws.Range("I132").Value = 2587.4546
$ws.Range("K132").Value = 7762.3638
$ws.Range("M132").Value = -5232.3638

# ---- Sheet: LTW (sheet7.xml) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6165.6665
$ws.Range("I7").Value = 5944.625
$ws.Range("K7").Value = 5944.625
$ws.Range("M7").Value = -5832.625
$ws.Range("H126").Value = 6165.6665
$ws.Range("I126").Value = 5944.625
$ws.Range("K126").Value = 17833.875
$ws.Range("M126").Value = -15363.875

# ---- Sheet: WVR (sheet8.xml) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2801
$ws.Range("I132").Value = 2702
$ws.Range("K132").Value = 8106
$ws.Range("M132").Value = -5576
$ws.Range("H141").Value = 299997.5
$ws.Range("J141").Value = 99995
$ws.Range("L141").Value = 99995
$ws.Range("N141").Value = -110355
